# Update countries & provincias Spain
# - Reorders "Costa de Marfil" so it sits right after "Senegal" (i.e. just
#   before "Cuba") in the country list, shifting Cuba/Islandia/Estonia down
#   one row and giving "Costa de Marfil" newly-updated statistics.
# - Refreshes the "datos actualizados" timestamp string.
# - Refreshes case/death/recovery counters for several countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 01:35"

# --- Estados Unidos (row 4) -------------------------------------------------
$ws.Cells.Item(4, 2).Value = 1408039
$ws.Cells.Item(4, 3).Value = 22205
$ws.Cells.Item(4, 4).Value = 293358
$ws.Cells.Item(4, 5).Value = 1031315
$ws.Cells.Item(4, 6).Value = 16473
$ws.Cells.Item(4, 7).Value = 1571
$ws.Cells.Item(4, 8).Value = 83366

# --- Brasil (row 10) ---------------------------------------------------
$ws.Cells.Item(10, 2).Value = 177602
$ws.Cells.Item(10, 3).Value = 8459
$ws.Cells.Item(10, 4).Value = 72597
$ws.Cells.Item(10, 5).Value = 92601
$ws.Cells.Item(10, 6).Value = 8318
$ws.Cells.Item(10, 7).Value = 779
$ws.Cells.Item(10, 8).Value = 12404

# --- India (row 15) ------------------------------------------------------
$ws.Cells.Item(15, 2).Value = 74292
$ws.Cells.Item(15, 3).Value = 3524
$ws.Cells.Item(15, 4).Value = 24420
$ws.Cells.Item(15, 5).Value = 47457
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 121
$ws.Cells.Item(15, 8).Value = 2415

# --- Nigeria (row 64) ------------------------------------------------------
$ws.Cells.Item(64, 2).Value = 4787
$ws.Cells.Item(64, 3).Value = 146
$ws.Cells.Item(64, 4).Value = 959
$ws.Cells.Item(64, 5).Value = 3670
$ws.Cells.Item(64, 6).Value = 7
$ws.Cells.Item(64, 7).Value = 8
$ws.Cells.Item(64, 8).Value = 158

# --- Reorder Costa de Marfil ahead of Cuba/Islandia/Estonia (rows 82-85) ---
# Row 82 becomes "Costa de Marfil" with freshly updated figures; the rows
# that used to hold Cuba/Islandia/Estonia data simply shift down by one,
# keeping all of their original values.
$ws.Cells.Item(82, 1).Value = "Costa de Marfil"
$ws.Cells.Item(82, 2).Value = 1857
$ws.Cells.Item(82, 3).Value = 127
$ws.Cells.Item(82, 4).Value = 820
$ws.Cells.Item(82, 5).Value = 1016
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 21

$ws.Cells.Item(83, 1).Value = "Cuba"
$ws.Cells.Item(83, 2).Value = 1804
$ws.Cells.Item(83, 3).Value = 21
$ws.Cells.Item(83, 4).Value = 1277
$ws.Cells.Item(83, 5).Value = 449
$ws.Cells.Item(83, 6).Value = 4
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = 78

$ws.Cells.Item(84, 1).Value = "Islandia"
$ws.Cells.Item(84, 2).Value = 1801
$ws.Cells.Item(84, 3).Value = 0
$ws.Cells.Item(84, 4).Value = 1776
$ws.Cells.Item(84, 5).Value = 15
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 10

$ws.Cells.Item(85, 1).Value = "Estonia"
$ws.Cells.Item(85, 2).Value = 1746
$ws.Cells.Item(85, 3).Value = 5
$ws.Cells.Item(85, 4).Value = 777
$ws.Cells.Item(85, 5).Value = 908
$ws.Cells.Item(85, 6).Value = 5
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 61

# --- Venezuela (row 129) --------------------------------------------------
$ws.Cells.Item(129, 2).Value = 423
$ws.Cells.Item(129, 3).Value = 1
$ws.Cells.Item(129, 4).Value = 220
$ws.Cells.Item(129, 5).Value = 193

# --- Polinesia Francesa (row 173) ------------------------------------------
$ws.Cells.Item(173, 4).Value = 58
$ws.Cells.Item(173, 5).Value = 2
